# CHIPS Budget.xlsx update — "Added Spice and Updated budget"
#
# Semantic changes applied:
#   1. Row 17 (Sensors section): item renamed "Coax Connector" -> "BNC"
#   2. Row 15 (Op Amps, Sensors section): Count 32 -> 16, Price/unit added (0.964)
#      -> its Total (G15) and the section Total (G23) and the grand Total (K4)
#         recalculate automatically from the formulas already in the sheet.
#   3. Row 26 (Microcontroller, Master Board section): annotate column H with "Sample"
#   4. Scroll/selection of the sheet view updated to reflect where the user was
#      working (top-left visible cell around row 19, active selection F38).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename item in row 17 from "Coax Connector" to "BNC"
$ws.Range("B17").Value = "BNC"

# 2. Update row 15 (Op Amps): reduce count and add a unit price
$ws.Range("E15").Value = 16
$ws.Range("F15").Value = 0.964

# 3. Add a "Sample" note next to the Microcontroller row
$ws.Range("H26").Value = "Sample"

# 4. Update the active window scroll position / selection
$win = $excel.ActiveWindow
$win.ScrollRow = 19
$win.ScrollColumn = 1
$ws.Range("F38").Select()
